$d = $word.ActiveDocument
$d.Paragraphs(2).Range.Delete()
